# Update countries & provincias Spain
# Applies updated COVID-19 case counts for several countries, and updates
# "Consejo Danes para los Refugiados" which now has more total cases than
# "Sudan" and therefore swaps rank (row) position with it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4: Estados Unidos ---
$ws.Range("B4").Value = 1188421
$ws.Range("C4").Value = 299
$ws.Range("E4").Value = 941225
$ws.Range("G4").Value = 4
$ws.Range("H4").Value = 68602

# --- Row 18: India ---
$ws.Range("B18").Value = 42670
$ws.Range("C18").Value = 165
$ws.Range("D18").Value = 11782
$ws.Range("E18").Value = 29493
$ws.Range("G18").Value = 4
$ws.Range("H18").Value = 1395

# --- Row 48: Colombia ---
$ws.Range("F48").Value = 122

# --- Row 80: Bulgaria ---
$ws.Range("B80").Value = 1632
$ws.Range("C80").Value = 14
$ws.Range("D80").Value = 321
$ws.Range("F80").Value = 40
$ws.Range("G80").Value = 1
$ws.Range("H80").Value = 74

# --- Rows 106/107: Consejo Danes para los Refugiados overtakes Sudan in the ranking ---
# Row 106 becomes "Consejo Danes para los Refugiados" with fresh updated data
$ws.Range("A106").Value = "Consejo Danes para los Refugiados"
$ws.Range("B106").Value = 682
$ws.Range("C106").Value = 8
$ws.Range("D106").Value = 80
$ws.Range("E106").Value = 568
$ws.Range("G106").Value = 1
$ws.Range("H106").Value = 34

# Row 107 becomes "Sudan" carrying the data it had while it was in row 106
$ws.Range("A107").Value = "Sudan"
$ws.Range("B107").Value = 678
$ws.Range("C107").Value = 86
$ws.Range("D107").Value = 61
$ws.Range("E107").Value = 576
$ws.Range("H107").Value = 41

# --- Row 121: Taiwan ---
$ws.Range("B121").Value = 437
$ws.Range("C121").Value = 5
$ws.Range("D121").Value = 334
$ws.Range("E121").Value = 97
